$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44495
$ws.Range("D3").Value = 44483
$ws.Range("J3").Value = 120
$ws.Range("D4").Value = 44477
$ws.Range("H4").Value = 'Sin especificar'
$ws.Range("I4").Value = 'Primera'
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 800
$ws.Range("L4").Value = 800
$ws.Range("M4").Value = 800
$ws.Range("N4").Value = '$/kilo (volumen en unidades)'
$ws.Range("O4").Value = 'Perú'
$ws.Range("P4").Value = 800
$ws.Range("D5").Value = 44497
$ws.Range("H5").Value = 'Sin especificar'
$ws.Range("J5").Value = 250
$ws.Range("K5").Value = 800
$ws.Range("L5").Value = 800
$ws.Range("M5").Value = 800
$ws.Range("N5").Value = '$/kilo (volumen en unidades)'
$ws.Range("O5").Value = 'Perú'
$ws.Range("P5").Value = 800
$ws.Range("D6").Value = 44305
$ws.Range("H6").Value = 'Sin especificar'
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 2500
$ws.Range("L6").Value = 2500
$ws.Range("M6").Value = 2500
$ws.Range("O6").Value = 'Perú'
$ws.Range("P6").Value = 2500
$ws.Range("D7").Value = 44312
$ws.Range("H7").Value = 'Sin especificar'
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 180
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 2500
$ws.Range("M7").Value = 2500
$ws.Range("O7").Value = 'Perú'
$ws.Range("P7").Value = 2500
$ws.Range("D8").Value = 44488
$ws.Range("J8").Value = 150
$ws.Range("D9").Value = 44510
$ws.Range("J9").Value = 250
$ws.Range("D10").Value = 44491
$ws.Range("J10").Value = 150
$ws.Range("D11").Value = 44167
$ws.Range("J11").Value = 400
$ws.Range("K11").Value = 5000
$ws.Range("L11").Value = 5000
$ws.Range("M11").Value = 5000
$ws.Range("O11").Value = 'Región de O''Higgins'
$ws.Range("P11").Value = 5000
$ws.Range("I12").Value = 'Segunda'
$ws.Range("J12").Value = 560
$ws.Range("K12").Value = 3000
$ws.Range("L12").Value = 3000
$ws.Range("M12").Value = 3000
$ws.Range("P12").Value = 3000
$ws.Range("I13").Value = 'Tercera'
$ws.Range("J13").Value = 450
$ws.Range("K13").Value = 2000
$ws.Range("L13").Value = 2000
$ws.Range("M13").Value = 2000
$ws.Range("P13").Value = 2000
$ws.Range("D14").Value = 44194
$ws.Range("I14").Value = 'Extra'
$ws.Range("J14").Value = 120
$ws.Range("K14").Value = 3500
$ws.Range("L14").Value = 3500
$ws.Range("M14").Value = 3500
$ws.Range("P14").Value = 3500
$ws.Range("D15").Value = 44194
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 3000
$ws.Range("L15").Value = 3000
$ws.Range("M15").Value = 3000
$ws.Range("O15").Value = 'Región de O''Higgins'
$ws.Range("P15").Value = 3000
$ws.Range("D16").Value = 44217
$ws.Range("I16").Value = 'Extra'
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 2500
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = 2500
$ws.Range("N16").Value = '$/unidad'
$ws.Range("O16").Value = 'Región de O''Higgins'
$ws.Range("P16").Value = 2500
$ws.Range("D17").Value = 44217
$ws.Range("J17").Value = 280
$ws.Range("K17").Value = 2000
$ws.Range("L17").Value = 2000
$ws.Range("M17").Value = 2000
$ws.Range("N17").Value = '$/unidad'
$ws.Range("O17").Value = 'Región de O''Higgins'
$ws.Range("P17").Value = 2000
$ws.Range("D18").Value = 44504
$ws.Range("D19").Value = 44223
$ws.Range("H19").Value = 'Americana O Klondike'
$ws.Range("J19").Value = 340
$ws.Range("D20").Value = 44223
$ws.Range("H20").Value = 'Americana O Klondike'
$ws.Range("J20").Value = 400
$ws.Range("D21").Value = 44223
$ws.Range("H21").Value = 'Americana O Klondike'
$ws.Range("I21").Value = 'Segunda'
$ws.Range("J21").Value = 300
$ws.Range("K21").Value = 1500
$ws.Range("L21").Value = 1500
$ws.Range("M21").Value = 1500
$ws.Range("P21").Value = 1500
$ws.Range("D22").Value = 44223
$ws.Range("H22").Value = 'Americana O Klondike'
$ws.Range("I22").Value = 'Tercera'
$ws.Range("J22").Value = 160
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = 1000
$ws.Range("P22").Value = 1000
